$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new AutoSklearn test-result row at the first empty row (row 12).
# Writing plain values/formulas here lets the cells inherit the column's
# default style (s="1" for C/F, s="4" for D/E/G/H) rather than copying the
# header row's style, which is what Insert() would do.
$ws.Range("A12").Value = "AutoSklearn - 60 seconds - 4 cores - 3GB RAM"
$ws.Range("B12").Value = $false
$ws.Range("C12").Value = 0.99166666666666603
$ws.Range("D12").Formula = "=1*60+10"
$ws.Range("E12").Value = 72
$ws.Range("F12").Value = 0.66518518518518499
$ws.Range("G12").Formula = "=1*60+5"
$ws.Range("H12").Value = 8

# Re-sort the whole results table (A2:I12) ascending by the Model column (A),
# matching the workbook's existing sortState so the new row lands
# alphabetically among the rest.
$dataRange = $ws.Range("A2:I12")
$keyColumn = $ws.Range("A2:A12")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyColumn)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Match the saved selection from the authored edit.
$ws.Range("H24").Select()
